$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("557:557").Insert()

$ws.Range("A557").Value = 3
$ws.Range("B557").Value = "Femacal de La Calera"
$ws.Range("C557").Value = "Coquimbo"
$ws.Range("D557").Value = 44984
$ws.Range("E557").Value = 5
$ws.Range("F557").Value = 100112027
$ws.Range("G557").Value = "Melón"
$ws.Range("H557").Value = "Tuna"
$ws.Range("I557").Value = "Primera"
$ws.Range("J557").Value = 550
$ws.Range("K557").Value = 1500
$ws.Range("L557").Value = 1500
$ws.Range("M557").Value = 1500
$ws.Range("N557").Value = "$/unidad"
$ws.Range("O557").Value = "Región de O'Higgins"
$ws.Range("P557").Value = 1500
$ws.Range("Q557").Value = 1
$ws.Range("R557").Value = "Hortaliza"
